# Auto-generated edit script applying the scheduled-runner price/profit refresh
# across the ALC, ARM, CRP, CUL, GSM, and LTW sheets (per upstream diff).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4814.8
$ws.Range("I74").Value = 4283
$ws.Range("K74").Value = 4283
$ws.Range("M74").Value = -3347
$ws.Range("H77").Value = 4814.8
$ws.Range("I77").Value = 4283
$ws.Range("K77").Value = 21415
$ws.Range("M77").Value = -16735

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 3447.7778
$ws.Range("I31").Value = 3447.7778
$ws.Range("K31").Value = 3447.7778
$ws.Range("M31").Value = -3153.7778
$ws.Range("H45").Value = 3201
$ws.Range("I45").Value = 2251.4167
$ws.Range("K45").Value = 2251.4167
$ws.Range("M45").Value = -1874.4167
$ws.Range("H74").Value = 29241.732
$ws.Range("I74").Value = 2303.5715
$ws.Range("J74").Value = 52812.625
$ws.Range("K74").Value = 2303.5715
$ws.Range("L74").Value = 52812.625
$ws.Range("M74").Value = -1429.5715
$ws.Range("N74").Value = -54560.625
$ws.Range("H77").Value = 29241.732
$ws.Range("I77").Value = 2303.5715
$ws.Range("J77").Value = 52812.625
$ws.Range("K77").Value = 11517.8575
$ws.Range("L77").Value = 264063.125
$ws.Range("M77").Value = -7149.8575
$ws.Range("N77").Value = -272799.125
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1467.8889
$ws.Range("I22").Value = 982.2222
$ws.Range("J22").Value = 1953.5555
$ws.Range("K22").Value = 982.2222
$ws.Range("L22").Value = 1953.5555
$ws.Range("M22").Value = -632.2222
$ws.Range("N22").Value = -2653.5555
$ws.Range("H31").Value = 188651.17
$ws.Range("I31").Value = 402093.2
$ws.Range("J31").Value = 36192.57
$ws.Range("K31").Value = 402093.2
$ws.Range("L31").Value = 36192.57
$ws.Range("M31").Value = -401798.2
$ws.Range("N31").Value = -36782.57
$ws.Range("H34").Value = 188651.17
$ws.Range("I34").Value = 402093.2
$ws.Range("J34").Value = 36192.57
$ws.Range("K34").Value = 402093.2
$ws.Range("L34").Value = 36192.57
$ws.Range("M34").Value = -401891.2
$ws.Range("N34").Value = -36596.57
$ws.Range("H99").Value = 10755.286
$ws.Range("I99").Value = 4547.8335
$ws.Range("K99").Value = 4547.8335
$ws.Range("M99").Value = -3049.8335
$ws.Range("H126").Value = 10755.286
$ws.Range("I126").Value = 4547.8335
$ws.Range("K126").Value = 13643.5005
$ws.Range("M126").Value = -11173.5005
$ws.Range("H132").Value = 102945380
$ws.Range("I132").Value = 4666.5
$ws.Range("J132").Value = 257356450
$ws.Range("K132").Value = 13999.5
$ws.Range("L132").Value = 772069350
$ws.Range("M132").Value = -11469.5
$ws.Range("N132").Value = -772074410
$ws.Range("H134").Value = 27032286
$ws.Range("I134").Value = 1723.9166
$ws.Range("K134").Value = 5171.7498
$ws.Range("M134").Value = -2636.7498

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 271
$ws.Range("I23").Value = 123
$ws.Range("J23").Value = 345
$ws.Range("K23").Value = 369
$ws.Range("L23").Value = 1035
$ws.Range("M23").Value = -134
$ws.Range("N23").Value = -1505
$ws.Range("H46").Value = 306.57144
$ws.Range("I46").Value = 309.2
$ws.Range("J46").Value = 300
$ws.Range("K46").Value = 927.5999999999999
$ws.Range("L46").Value = 900
$ws.Range("M46").Value = -836.5999999999999
$ws.Range("N46").Value = -1082
$ws.Range("H110").Value = 16172.25
$ws.Range("I110").Value = 1889
$ws.Range("J110").Value = 30455.5
$ws.Range("K110").Value = 5667
$ws.Range("L110").Value = 91366.5
$ws.Range("M110").Value = -1577
$ws.Range("N110").Value = -99546.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12689.2
$ws.Range("I70").Value = 12850
$ws.Range("J70").Value = 12582
$ws.Range("K70").Value = 12850
$ws.Range("L70").Value = 12582
$ws.Range("M70").Value = -12580
$ws.Range("N70").Value = -13122
$ws.Range("H73").Value = 12689.2
$ws.Range("I73").Value = 12850
$ws.Range("J73").Value = 12582
$ws.Range("K73").Value = 12850
$ws.Range("L73").Value = 12582
$ws.Range("M73").Value = -11914
$ws.Range("N73").Value = -14454
$ws.Range("H80").Value = 4766
$ws.Range("I80").Value = 5207.5
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 5207.5
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -4209.5
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 4766
$ws.Range("I83").Value = 5207.5
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 26037.5
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -21045.5
$ws.Range("N83").Value = -24984
$ws.Range("H122").Value = 2129.4167
$ws.Range("I122").Value = 1790.8096
$ws.Range("J122").Value = 4499.6665
$ws.Range("K122").Value = 5372.4288
$ws.Range("L122").Value = 13498.9995
$ws.Range("M122").Value = -2922.4288
$ws.Range("N122").Value = -18398.9995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 988.8461
$ws.Range("I22").Value = 664.0909
$ws.Range("J22").Value = 1227
$ws.Range("K22").Value = 664.0909
$ws.Range("L22").Value = 1227
$ws.Range("M22").Value = -369.0909
$ws.Range("N22").Value = -1817
$ws.Range("H25").Value = 509998.84
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 509998.84
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 509998.84
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -510458.84
$ws.Range("H27").Value = 988.8461
$ws.Range("I27").Value = 664.0909
$ws.Range("J27").Value = 1227
$ws.Range("K27").Value = 664.0909
$ws.Range("L27").Value = 1227
$ws.Range("M27").Value = -557.0909
$ws.Range("N27").Value = -1441
$ws.Range("H40").Value = 4734.643
$ws.Range("I40").Value = 1716.3334
$ws.Range("J40").Value = 6998.375
$ws.Range("K40").Value = 1716.3334
$ws.Range("L40").Value = 6998.375
$ws.Range("M40").Value = -1580.3334
$ws.Range("N40").Value = -7270.375
$ws.Range("H46").Value = 3024.4375
$ws.Range("I46").Value = 1400
$ws.Range("J46").Value = 3399.3076
$ws.Range("K46").Value = 1400
$ws.Range("L46").Value = 3399.3076
$ws.Range("M46").Value = -1212
$ws.Range("N46").Value = -3775.3076
$ws.Range("H56").Value = 10524.286
$ws.Range("I56").Value = 8666.666999999999
$ws.Range("K56").Value = 8666.666999999999
$ws.Range("M56").Value = -7975.666999999999
$ws.Range("H99").Value = 29369.166
$ws.Range("I99").Value = 29369.166
$ws.Range("K99").Value = 29369.166
$ws.Range("M99").Value = -26374.166
$ws.Range("H122").Value = 12272.272
$ws.Range("I122").Value = 15000
$ws.Range("K122").Value = 45000
$ws.Range("M122").Value = -42550
